$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C values (rows 2-12)
$ws.Range("C2").Value = -8.804082051850855
$ws.Range("C3").Value = 0.001461578300222754
$ws.Range("C4").Value = -2.506954065174796
$ws.Range("C5").Value = -1.533333595940348
$ws.Range("C6").Value = -0.0008846351120155305
$ws.Range("C7").Value = 0.7370270550018176
$ws.Range("C8").Value = -2.130658382258844
$ws.Range("C9").Value = -1.571411290919059
$ws.Range("C10").Value = -0.8815430758986622
$ws.Range("C11").Value = -0.02860128806787543
$ws.Range("C12").Value = -1.338265054859221

# Column B values (rows 13-22)
$ws.Range("B13").Value = -0.7292374349199235
$ws.Range("B14").Value = -71.11691986769438
$ws.Range("B15").Value = -2.293386437115259
$ws.Range("B16").Value = -0.4560600868171605
$ws.Range("B17").Value = -0.345623665722087
$ws.Range("B18").Value = 0.1470122990431264
$ws.Range("B19").Value = -0.4087767791788792
$ws.Range("B20").Value = -0.4447727375663817
$ws.Range("B21").Value = -0.02860128808242735
$ws.Range("B22").Value = -1.237662993371487
